$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header + values for the new "Locator Type" column
$ws.Range("D1").Value = "Locator Type"
$ws.Range("D2").Value = "CSS"
$ws.Range("D3").Value = "CSS"
$ws.Range("D4").Value = "CSS"
$ws.Range("D5").Value = "CSS"
$ws.Range("D6").Value = "CSS"
$ws.Range("D7").Value = "CSS"
$ws.Range("D8").Value = "CSS"
$ws.Range("D9").Value = "Xpath"
$ws.Range("D10").Value = "Xpath"
$ws.Range("D11").Value = "CSS"
$ws.Range("D12").Value = "CSS"
$ws.Range("D13").Value = "CSS"
$ws.Range("D14").Value = "Xpath"

# Match column D's cell formatting to the rest of the plain data cells
# (text number format, default font) used throughout the sheet.
$ws.Range("D1:D14").NumberFormat = "@"

# Set the new column's width to fit its content ("Locator Type" header).
$ws.Columns.Item(4).ColumnWidth = 9.666666666666666

# Update view state to match the post-edit selection/scroll position.
$ws.Range("D14").Select()
